# Auto-applied edit script for horarios-141 workbook update
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 20:00:07"
$ws.Cells.Item(3, 1).Value = "Total filas: 311"
$ws.Cells.Item(15, 3).Value = "225_GOMEZ"
$ws.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws.Cells.Item(23, 1).Value = "06:46:40"
$ws.Cells.Item(23, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(23, 4).Value = 35
$ws.Cells.Item(24, 1).Value = "06:15:23"
$ws.Cells.Item(24, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(24, 4).Value = 66
$ws.Cells.Item(55, 1).Value = "08:14:55"
$ws.Cells.Item(55, 3).Value = "215B_EL PATO"
$ws.Cells.Item(55, 4).Value = 39
$ws.Cells.Item(56, 1).Value = "08:49:06"
$ws.Cells.Item(56, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(56, 4).Value = 4
$ws.Cells.Item(64, 1).Value = "08:57:42"
$ws.Cells.Item(64, 3).Value = "14_ABASTO"
$ws.Cells.Item(64, 4).Value = 20
$ws.Cells.Item(65, 1).Value = "08:49:06"
$ws.Cells.Item(65, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(65, 4).Value = 28
$ws.Cells.Item(66, 1).Value = "08:57:42"
$ws.Cells.Item(66, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(66, 4).Value = 21
$ws.Cells.Item(67, 1).Value = "08:49:06"
$ws.Cells.Item(67, 3).Value = "14_ABASTO"
$ws.Cells.Item(67, 4).Value = 29
$ws.Cells.Item(71, 1).Value = "08:14:55"
$ws.Cells.Item(71, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(71, 4).Value = 77
$ws.Cells.Item(72, 1).Value = "08:49:06"
$ws.Cells.Item(72, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(72, 4).Value = 42
$ws.Cells.Item(90, 3).Value = "14_ABASTO"
$ws.Cells.Item(91, 3).Value = "15_ABASTO"
$ws.Cells.Item(140, 1).Value = "12:18:38"
$ws.Cells.Item(140, 3).Value = "215C_EL PATO"
$ws.Cells.Item(140, 4).Value = 45
$ws.Cells.Item(141, 1).Value = "12:43:13"
$ws.Cells.Item(141, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(141, 4).Value = 20
$ws.Cells.Item(148, 1).Value = "12:58:23"
$ws.Cells.Item(148, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(148, 4).Value = 27
$ws.Cells.Item(149, 1).Value = "12:43:13"
$ws.Cells.Item(149, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(149, 4).Value = 42
$ws.Cells.Item(172, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(173, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(191, 3).Value = "10_OLMOS"
$ws.Cells.Item(192, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(199, 1).Value = "14:58:38"
$ws.Cells.Item(199, 3).Value = "14_ABASTO"
$ws.Cells.Item(199, 4).Value = 67
$ws.Cells.Item(200, 1).Value = "16:02:30"
$ws.Cells.Item(200, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(200, 4).Value = 3
$ws.Cells.Item(223, 1).Value = "16:34:05"
$ws.Cells.Item(223, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(223, 4).Value = 33
$ws.Cells.Item(224, 1).Value = "16:57:38"
$ws.Cells.Item(224, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(224, 4).Value = 10
$ws.Cells.Item(249, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(250, 3).Value = "15_ABASTO"
$ws.Cells.Item(280, 1).Value = "18:01:05"
$ws.Cells.Item(280, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(280, 4).Value = 82
$ws.Cells.Item(281, 1).Value = "19:14:15"
$ws.Cells.Item(281, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(281, 4).Value = 9
$ws.Cells.Item(291, 1).Value = "20:00:07"
$ws.Cells.Item(291, 4).Value = 7
$ws.Cells.Item(292, 1).Value = "20:00:07"
$ws.Cells.Item(292, 2).Value = "20:08"
$ws.Cells.Item(292, 4).Value = 8
$ws.Cells.Item(293, 1).Value = "19:45:00"
$ws.Cells.Item(293, 2).Value = "20:09"
$ws.Cells.Item(293, 4).Value = 24
$ws.Cells.Item(294, 1).Value = "19:14:15"
$ws.Cells.Item(294, 2).Value = "20:11"
$ws.Cells.Item(294, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(294, 4).Value = 57
$ws.Cells.Item(295, 1).Value = "18:50:27"
$ws.Cells.Item(295, 2).Value = "20:12"
$ws.Cells.Item(295, 4).Value = 82
$ws.Cells.Item(296, 1).Value = "20:00:07"
$ws.Cells.Item(296, 2).Value = "20:13"
$ws.Cells.Item(296, 3).Value = "14_ABASTO"
$ws.Cells.Item(296, 4).Value = 13
$ws.Cells.Item(297, 1).Value = "20:00:07"
$ws.Cells.Item(297, 2).Value = "20:22"
$ws.Cells.Item(297, 3).Value = "15_ABASTO"
$ws.Cells.Item(297, 4).Value = 22
$ws.Cells.Item(298, 1).Value = "18:50:27"
$ws.Cells.Item(298, 2).Value = "20:30"
$ws.Cells.Item(298, 4).Value = 100
$ws.Cells.Item(299, 1).Value = "20:00:07"
$ws.Cells.Item(299, 2).Value = "20:31"
$ws.Cells.Item(299, 3).Value = "10_OLMOS"
$ws.Cells.Item(299, 4).Value = 31
$ws.Cells.Item(300, 1).Value = "20:00:07"
$ws.Cells.Item(300, 2).Value = "20:34"
$ws.Cells.Item(300, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(300, 4).Value = 34
$ws.Cells.Item(301, 1).Value = "20:00:07"
$ws.Cells.Item(301, 2).Value = "20:46"
$ws.Cells.Item(301, 3).Value = "17_ROMERO"
$ws.Cells.Item(301, 4).Value = 46
$ws.Cells.Item(302, 1).Value = "18:50:27"
$ws.Cells.Item(302, 2).Value = "20:47"
$ws.Cells.Item(302, 4).Value = 117
$ws.Cells.Item(303, 1).Value = "20:00:07"
$ws.Cells.Item(303, 2).Value = "20:48"
$ws.Cells.Item(303, 3).Value = "215B_EL PATO"
$ws.Cells.Item(303, 4).Value = 48
$ws.Cells.Item(304, 1).Value = "19:14:15"
$ws.Cells.Item(304, 2).Value = "20:50"
$ws.Cells.Item(304, 4).Value = 96
$ws.Cells.Item(305, 2).Value = "20:52"
$ws.Cells.Item(305, 3).Value = "17_ROMERO"
$ws.Cells.Item(305, 4).Value = 67
$ws.Cells.Item(306, 2).Value = "20:55"
$ws.Cells.Item(306, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(306, 4).Value = 70
$ws.Cells.Item(307, 1).Value = "19:45:00"
$ws.Cells.Item(307, 2).Value = "20:56"
$ws.Cells.Item(307, 4).Value = 71
$ws.Cells.Item(308, 1).Value = "20:00:07"
$ws.Cells.Item(308, 2).Value = "20:57"
$ws.Cells.Item(308, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(308, 4).Value = 57
$ws.Cells.Item(309, 1).Value = "20:00:07"
$ws.Cells.Item(309, 2).Value = "21:07"
$ws.Cells.Item(309, 3).Value = "10_OLMOS"
$ws.Cells.Item(309, 4).Value = 67
$ws.Cells.Item(310, 1).Value = "20:00:07"
$ws.Cells.Item(310, 2).Value = "21:10"
$ws.Cells.Item(310, 3).Value = "15_ABASTO"
$ws.Cells.Item(310, 4).Value = 70
$ws.Cells.Item(311, 2).Value = "21:28"
$ws.Cells.Item(311, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(311, 4).Value = 103
$ws.Cells.Item(312, 1).Value = "20:00:07"
$ws.Cells.Item(312, 2).Value = "21:29"
$ws.Cells.Item(312, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(312, 4).Value = 89
$ws.Cells.Item(312, 5).Value = "LP1912"
$ws.Cells.Item(313, 1).Value = "20:00:07"
$ws.Cells.Item(313, 2).Value = "21:33"
$ws.Cells.Item(313, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(313, 4).Value = 93
$ws.Cells.Item(313, 5).Value = "LP1912"
$ws.Cells.Item(314, 1).Value = "20:00:07"
$ws.Cells.Item(314, 2).Value = "21:34"
$ws.Cells.Item(314, 3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws.Cells.Item(314, 4).Value = 94
$ws.Cells.Item(314, 5).Value = "LP1912"
$ws.Cells.Item(315, 1).Value = "20:00:07"
$ws.Cells.Item(315, 2).Value = "21:46"
$ws.Cells.Item(315, 3).Value = "14X44_ABASTO"
$ws.Cells.Item(315, 4).Value = 106
$ws.Cells.Item(315, 5).Value = "LP1912"
$ws.Cells.Item(316, 1).Value = "20:00:07"
$ws.Cells.Item(316, 2).Value = "21:48"
$ws.Cells.Item(316, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(316, 4).Value = 108
$ws.Cells.Item(316, 5).Value = "LP1912"

$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 20:00:07"
$ws.Cells.Item(51, 1).Value = "20:00:07"
$ws.Cells.Item(51, 4).Value = 7
$ws.Cells.Item(53, 1).Value = "20:00:07"
$ws.Cells.Item(53, 4).Value = 48

$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 20:00:07"
$ws.Cells.Item(3, 1).Value = "Total filas: 42"
$ws.Cells.Item(46, 1).Value = "20:00:07"
$ws.Cells.Item(46, 2).Value = "20:00"
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(47, 1).Value = "20:00:07"
$ws.Cells.Item(47, 2).Value = "20:52"
$ws.Cells.Item(47, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(47, 4).Value = 52
$ws.Cells.Item(47, 5).Value = "L6203"
